# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the "In Translation"
# status became "Ready for handoff" (wherever it appears) and the
# "Latest HO/Handoff" timestamps were bumped forward a few dozen seconds.
# Excel auto-sized the Status/zh-cn/de-de columns to fit the new, longer
# "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-17 18:38:25"
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-17 18:38:20"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.3

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-17 18:38:25"
$dede.Range("C1").EntireColumn.ColumnWidth = 16.3
